$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K" (formerly derived from "Strike#"). Regenerated
# values (s_vals) replace the previous data for rows 2-9.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 1
